$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 data
$ws.Range("A16").Value = 45958
$ws.Range("B16").Value = 5603
$ws.Range("C16").Value = 4351
$ws.Range("D16").Value = 4029
$ws.Range("E16").Value = 236
$ws.Range("F16").Value = 49
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 2

# Row 17 data
$ws.Range("A17").Value = 45959
$ws.Range("B17").Value = 5602
$ws.Range("C17").Value = 4355
$ws.Range("D17").Value = 4026
$ws.Range("E17").Value = 238
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 1

# Update the active selection to match the diff
$ws.Range("A17:I17").Select()
